# daily auto push: 2026-02-06 19:05 UTC
# Insert a new data row for 2026/02/07 (Sat) just above the 2026/12/29 row,
# shifting the existing rows 793-834 down to 794-835.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 793 down by one row.
$ws.Rows.Item(793).Insert()

# Fill in the freshly inserted row with the new record.
# Column A holds dates stored as plain text (matches the rest of the sheet),
# so force a Text format before assigning, then clear the format again so the
# cell ends up with no explicit style (same as its neighbours) while keeping
# the value as literal text instead of being auto-converted to a date serial.
$dateCell = $ws.Range("A793")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/02/07"
$dateCell.ClearFormats()

$ws.Range("B793").Value = "土"
$ws.Range("C793").Value = 1
$ws.Range("D793").Value = 201
